$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header - copy style from the neighboring header cell (G1) so it reuses
# the same cell format (bold, border, centered) instead of minting a new one
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for new "Save" column
$values = @(0, 0, 0, 1, 1, 0, 0, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
